# Replace the product image filenames from .jpg to .webp across all sheets
# (commit message: "prova foto webP file excel").

$wb = $excel.ActiveWorkbook

function Set-ImageCell {
    param($sheetName, $cellRef, $newValue)
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $newValue
}

# Antipasti
Set-ImageCell "Antipasti" "D5"  "pata_pig.webp"
Set-ImageCell "Antipasti" "D8"  "mangia_e_bevi_di_pizza.webp"
Set-ImageCell "Antipasti" "D9"  "sushi_pizza.webp"
Set-ImageCell "Antipasti" "D15" "arancinette_alla_carne.webp"

# Pizze Classiche
Set-ImageCell "Pizze Classiche" "D4"  "margherita.webp"
Set-ImageCell "Pizze Classiche" "D6"  "romana.webp"
Set-ImageCell "Pizze Classiche" "D8"  "capricciosa.webp"
Set-ImageCell "Pizze Classiche" "D10" "sfiziosa.webp"
Set-ImageCell "Pizze Classiche" "D11" "parmigiana.webp"
Set-ImageCell "Pizze Classiche" "D16" "american_chips.webp"

# Pizze Gourmet
Set-ImageCell "Pizze Gourmet" "D11" "bouquet.webp"
Set-ImageCell "Pizze Gourmet" "D13" "duchessa.webp"

# Hamburger
Set-ImageCell "Hamburger" "D2"  "gigi_burger.webp"
Set-ImageCell "Hamburger" "D5"  "porchettone.webp"
Set-ImageCell "Hamburger" "D10" "porchi_burger.webp"
Set-ImageCell "Hamburger" "D11" "burger_pig.webp"
Set-ImageCell "Hamburger" "D13" "cheese_burger.webp"

# Secondi
Set-ImageCell "Secondi" "D3" "pollo_panato_alla_griglia.webp"
Set-ImageCell "Secondi" "D4" "spiedini_al_pistacchio.webp"
Set-ImageCell "Secondi" "D5" "spiedini_siciliani.webp"
Set-ImageCell "Secondi" "D7" "tagliata_di_scottona.webp"

# Dolci
Set-ImageCell "Dolci" "D2" "tortino_cuore_caldo.webp"
Set-ImageCell "Dolci" "D3" "zeppola.webp"
Set-ImageCell "Dolci" "D4" "flauto.webp"

# Restore the active sheet / selection to match the saved workbook state
# (Pizze Classiche, cell D4 selected).
$wsActive = $wb.Worksheets.Item("Pizze Classiche")
$wsActive.Activate()
$wsActive.Range("D4").Select()
